$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Rename sheet from "shiccsd" to "dwfuy"
$ws.Name = "dwfuy"

# Update criterion values (row-by-row, per the diff)
# Row 1
$ws.Cells.Item(1, 10).Value = 32.26990509033203

# Row 2
$ws.Cells.Item(2, 2).Value = 1869
$ws.Cells.Item(2, 4).Value = 1864
$ws.Cells.Item(2, 5).Value = 4
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 100
$ws.Cells.Item(2, 8).Value = 99.78586723768737
$ws.Cells.Item(2, 9).Value = 0.002144772117962467
$ws.Cells.Item(2, 10).Value = 43.60518383979797

# Row 3
$ws.Cells.Item(3, 2).Value = 2083
$ws.Cells.Item(3, 10).Value = 39.43747496604919

# Row 4
$ws.Cells.Item(4, 2).Value = 2594
$ws.Cells.Item(4, 4).Value = 2566
$ws.Cells.Item(4, 5).Value = 27
$ws.Cells.Item(4, 6).Value = 5
$ws.Cells.Item(4, 7).Value = 99.80552314274601
$ws.Cells.Item(4, 8).Value = 98.95873505591979
$ws.Cells.Item(4, 9).Value = 0.01244167962674961
$ws.Cells.Item(4, 10).Value = 34.94165802001953

# Row 5
$ws.Cells.Item(5, 2).Value = 2016
$ws.Cells.Item(5, 4).Value = 2013
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 13
$ws.Cells.Item(5, 7).Value = 99.3583415597236
$ws.Cells.Item(5, 8).Value = 99.90074441687345
$ws.Cells.Item(5, 9).Value = 0.00740009866798224
$ws.Cells.Item(5, 10).Value = 42.01006698608398

# Row 6
$ws.Cells.Item(6, 2).Value = 1758
$ws.Cells.Item(6, 4).Value = 1749
$ws.Cells.Item(6, 5).Value = 8
$ws.Cells.Item(6, 6).Value = 13
$ws.Cells.Item(6, 7).Value = 99.2622020431328
$ws.Cells.Item(6, 8).Value = 99.54467842914057
$ws.Cells.Item(6, 9).Value = 0.01191151446398185
$ws.Cells.Item(6, 10).Value = 35.57951617240906

# Row 7
$ws.Cells.Item(7, 2).Value = 2532
$ws.Cells.Item(7, 4).Value = 2531
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 100
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 39.08774733543396

# Row 8
$ws.Cells.Item(8, 2).Value = 2124
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 8).Value = 99.95289684408856
$ws.Cells.Item(8, 9).Value = 0.0009416195856873823
$ws.Cells.Item(8, 10).Value = 33.75424337387085

# Row 9
$ws.Cells.Item(9, 10).Value = 37.7343442440033

# Row 10
$ws.Cells.Item(10, 2).Value = 1814
$ws.Cells.Item(10, 5).Value = 20
$ws.Cells.Item(10, 8).Value = 98.89685603971319
$ws.Cells.Item(10, 9).Value = 0.0116991643454039
$ws.Cells.Item(10, 10).Value = 40.09796524047852

# Row 11
$ws.Cells.Item(11, 2).Value = 1880
$ws.Cells.Item(11, 4).Value = 1876
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 2
$ws.Cells.Item(11, 7).Value = 99.89350372736955
$ws.Cells.Item(11, 8).Value = 99.84034060670569
$ws.Cells.Item(11, 9).Value = 0.002660989888238425
$ws.Cells.Item(11, 10).Value = 29.73125195503235

# Row 12
$ws.Cells.Item(12, 10).Value = 34.34687185287476

# Row 13
$ws.Cells.Item(13, 2).Value = 2389
$ws.Cells.Item(13, 4).Value = 2388
$ws.Cells.Item(13, 6).Value = 23
$ws.Cells.Item(13, 7).Value = 99.0460389879718
$ws.Cells.Item(13, 9).Value = 0.009535655058043118
$ws.Cells.Item(13, 10).Value = 33.43625259399414

# Row 14
$ws.Cells.Item(14, 10).Value = 32.32281112670898

# Row 15
$ws.Cells.Item(15, 10).Value = 34.04889917373657

# Row 16
$ws.Cells.Item(16, 2).Value = 1988
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 8).Value = 99.94967287367891
$ws.Cells.Item(16, 9).Value = 0.0005032712632108706
$ws.Cells.Item(16, 10).Value = 35.9206268787384

# Row 17
$ws.Cells.Item(17, 10).Value = 37.22344851493835

# Row 18
$ws.Cells.Item(18, 10).Value = 36.21926259994507

# Row 19
$ws.Cells.Item(19, 10).Value = 36.46384978294373

# Row 20
$ws.Cells.Item(20, 2).Value = 1613
$ws.Cells.Item(20, 4).Value = 1612
$ws.Cells.Item(20, 6).Value = 6
$ws.Cells.Item(20, 7).Value = 99.62917181705809
$ws.Cells.Item(20, 9).Value = 0.003705991352686844
$ws.Cells.Item(20, 10).Value = 28.86901164054871

# Row 21
$ws.Cells.Item(21, 2).Value = 2623
$ws.Cells.Item(21, 5).Value = 25
$ws.Cells.Item(21, 8).Value = 99.0465293668955
$ws.Cells.Item(21, 9).Value = 0.01076509034986544
$ws.Cells.Item(21, 10).Value = 35.9991626739502

# Row 22
$ws.Cells.Item(22, 2).Value = 1933
$ws.Cells.Item(22, 4).Value = 1932
$ws.Cells.Item(22, 6).Value = 30
$ws.Cells.Item(22, 7).Value = 98.47094801223241
$ws.Cells.Item(22, 9).Value = 0.01528273051451859
$ws.Cells.Item(22, 10).Value = 42.13510394096375

# Row 23
$ws.Cells.Item(23, 2).Value = 2062
$ws.Cells.Item(23, 4).Value = 2060
$ws.Cells.Item(23, 5).Value = 1
$ws.Cells.Item(23, 6).Value = 75
$ws.Cells.Item(23, 7).Value = 96.48711943793911
$ws.Cells.Item(23, 8).Value = 99.95147986414362
$ws.Cells.Item(23, 9).Value = 0.03558052434456929
$ws.Cells.Item(23, 10).Value = 33.41648864746094

# Row 24
$ws.Cells.Item(24, 2).Value = 3007
$ws.Cells.Item(24, 4).Value = 2955
$ws.Cells.Item(24, 5).Value = 51
$ws.Cells.Item(24, 6).Value = 24
$ws.Cells.Item(24, 7).Value = 99.19436052366567
$ws.Cells.Item(24, 8).Value = 98.30339321357286
$ws.Cells.Item(24, 9).Value = 0.02516778523489933
$ws.Cells.Item(24, 10).Value = 42.64461779594421

# Row 25
$ws.Cells.Item(25, 2).Value = 2648
$ws.Cells.Item(25, 4).Value = 2647
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 8
$ws.Cells.Item(25, 7).Value = 99.69868173258004
$ws.Cells.Item(25, 8).Value = 100
$ws.Cells.Item(25, 9).Value = 0.003012048192771084
$ws.Cells.Item(25, 10).Value = 36.2281768321991

# Row 26
$ws.Cells.Item(26, 2).Value = 1832
$ws.Cells.Item(26, 4).Value = 1827
$ws.Cells.Item(26, 5).Value = 4
$ws.Cells.Item(26, 6).Value = 32
$ws.Cells.Item(26, 7).Value = 98.27864443249058
$ws.Cells.Item(26, 8).Value = 99.7815401419989
$ws.Cells.Item(26, 9).Value = 0.01935483870967742
$ws.Cells.Item(26, 10).Value = 32.87388706207275

# Row 27
$ws.Cells.Item(27, 2).Value = 2947
$ws.Cells.Item(27, 4).Value = 2939
$ws.Cells.Item(27, 5).Value = 7
$ws.Cells.Item(27, 6).Value = 15
$ws.Cells.Item(27, 7).Value = 99.49221394719025
$ws.Cells.Item(27, 8).Value = 99.76238968092329
$ws.Cells.Item(27, 9).Value = 0.007445008460236886
$ws.Cells.Item(27, 10).Value = 38.88775682449341

# Row 28
$ws.Cells.Item(28, 2).Value = 3011
$ws.Cells.Item(28, 5).Value = 6
$ws.Cells.Item(28, 8).Value = 99.80066445182725
$ws.Cells.Item(28, 9).Value = 0.001996672212978369
$ws.Cells.Item(28, 10).Value = 36.25552558898926

# Row 29
$ws.Cells.Item(29, 2).Value = 2650
$ws.Cells.Item(29, 4).Value = 2635
$ws.Cells.Item(29, 5).Value = 14
$ws.Cells.Item(29, 6).Value = 14
$ws.Cells.Item(29, 7).Value = 99.4714986787467
$ws.Cells.Item(29, 8).Value = 99.4714986787467
$ws.Cells.Item(29, 9).Value = 0.01056603773584906
$ws.Cells.Item(29, 10).Value = 37.97565245628357

# Row 30
$ws.Cells.Item(30, 2).Value = 2753
$ws.Cells.Item(30, 5).Value = 5
$ws.Cells.Item(30, 8).Value = 99.81831395348837
$ws.Cells.Item(30, 9).Value = 0.001819505094614265
$ws.Cells.Item(30, 10).Value = 36.00034880638123

# Row 31
$ws.Cells.Item(31, 2).Value = 3249
$ws.Cells.Item(31, 4).Value = 3248
$ws.Cells.Item(31, 6).Value = 2
$ws.Cells.Item(31, 7).Value = 99.93846153846154
$ws.Cells.Item(31, 9).Value = 0.0006151953245155337
$ws.Cells.Item(31, 10).Value = 35.80213499069214

# Row 32
$ws.Cells.Item(32, 2).Value = 2264
$ws.Cells.Item(32, 4).Value = 2256
$ws.Cells.Item(32, 5).Value = 7
$ws.Cells.Item(32, 6).Value = 5
$ws.Cells.Item(32, 7).Value = 99.77885891198585
$ws.Cells.Item(32, 8).Value = 99.69067609368095
$ws.Cells.Item(32, 9).Value = 0.005305039787798408
$ws.Cells.Item(32, 10).Value = 39.08246159553528

# Row 33
$ws.Cells.Item(33, 2).Value = 3363
$ws.Cells.Item(33, 4).Value = 3362
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(33, 7).Value = 100
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 10).Value = 39.46951222419739

# Row 34
$ws.Cells.Item(34, 10).Value = 38.75755262374878

# Row 35
$ws.Cells.Item(35, 10).Value = 43.05935978889465

# Row 36
$ws.Cells.Item(36, 4).Value = 2420
$ws.Cells.Item(36, 5).Value = 3
$ws.Cells.Item(36, 6).Value = 6
$ws.Cells.Item(36, 7).Value = 99.75267930750206
$ws.Cells.Item(36, 8).Value = 99.87618654560463
$ws.Cells.Item(36, 9).Value = 0.003708281829419036
$ws.Cells.Item(36, 10).Value = 35.71985626220703

# Row 37
$ws.Cells.Item(37, 2).Value = 2473
$ws.Cells.Item(37, 4).Value = 2469
$ws.Cells.Item(37, 5).Value = 3
$ws.Cells.Item(37, 6).Value = 13
$ws.Cells.Item(37, 7).Value = 99.47622884770347
$ws.Cells.Item(37, 8).Value = 99.87864077669903
$ws.Cells.Item(37, 9).Value = 0.00644381796214257
$ws.Cells.Item(37, 10).Value = 39.78744673728943

# Row 38
$ws.Cells.Item(38, 2).Value = 2606
$ws.Cells.Item(38, 4).Value = 2604
$ws.Cells.Item(38, 5).Value = 1
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 100
$ws.Cells.Item(38, 8).Value = 99.9616122840691
$ws.Cells.Item(38, 9).Value = 0.0003838771593090211
$ws.Cells.Item(38, 10).Value = 30.80963039398193

# Row 39
$ws.Cells.Item(39, 2).Value = 2077
$ws.Cells.Item(39, 5).Value = 29
$ws.Cells.Item(39, 8).Value = 98.60308285163777
$ws.Cells.Item(39, 9).Value = 0.01656113005358013
$ws.Cells.Item(39, 10).Value = 39.27287030220032

# Row 40
$ws.Cells.Item(40, 10).Value = 40.54674196243286

# Row 41
$ws.Cells.Item(41, 10).Value = 36.33669281005859

# Row 42
$ws.Cells.Item(42, 2).Value = 1779
$ws.Cells.Item(42, 4).Value = 1778
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(42, 6).Value = 1
$ws.Cells.Item(42, 7).Value = 99.94378864530636
$ws.Cells.Item(42, 8).Value = 100
$ws.Cells.Item(42, 9).Value = 0.0005617977528089888
$ws.Cells.Item(42, 10).Value = 35.81988716125488

# Row 43
$ws.Cells.Item(43, 2).Value = 3079
$ws.Cells.Item(43, 4).Value = 3078
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(43, 7).Value = 100
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 43.74973726272583

# Row 44
$ws.Cells.Item(44, 10).Value = 37.19406127929688
